# Ethiopia_Section2_other_info: "Changed the localize handlebars function
# to expose the data model and the calculates" - the question prompts on
# the "survey" sheet now read {{data.name}} instead of {{name}}.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
[void]$survey.Activate()

$survey.Range("G2").Value  = 'What is {{data.name}}''s relationship to the household head?'
$survey.Range("G3").Value  = 'What is {{data.name}}''s sex?'
$survey.Range("G4").Value  = 'What is {{data.name}}''s age?'
$survey.Range("G6").Value  = 'Marital status of {{data.name}}.'
$survey.Range("G9").Value  = 'For how many months during the last 12 months was {{data.name}} away from the household?'
$survey.Range("G11").Value = 'In what region was {{data.name}} born?'
$survey.Range("G13").Value = 'What is {{data.name}}''s main religion?'

# Restore the per-sheet selection state recorded in the workbook (the
# cell that was active the last time each sheet was used in Excel).
[void]$survey.Range("G13").Select()

$choices = $wb.Worksheets.Item("choices")
[void]$choices.Activate()
[void]$choices.Range("C19").Select()

$settings = $wb.Worksheets.Item("settings")
[void]$settings.Activate()
[void]$settings.Range("A10").Select()

$model = $wb.Worksheets.Item("model")
[void]$model.Activate()
[void]$model.Range("H15").Select()

# "initial" remains the active tab, as in the original workbook.
$initial = $wb.Worksheets.Item("initial")
[void]$initial.Activate()
